$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.977869842962853, 1.0222959807518261, 0.6292826406311872, 2.64824563716967)"
$ws.Range("C2").Value = "JSU(-1.4382036252289248, 1.2022004731650933, 2.2972591131821076, 4.2882976805457265)"
$ws.Range("D2").Value = "JSU(-0.953943921961079, 1.1015774363053588, 0.7086876866351115, 2.5698407400848913)"
$ws.Range("E2").Value = "NCT(3.267177661873953, 1.5925321342597303, -0.018666485496264724, 4.749215643126205)"
